# Update column G ("K") values on Sheet1 for rows 2-13.
# These values were regenerated after switching from "Strike#" to "K"
# and recalculating std/mean (s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 0
    4  = 5
    5  = 4
    6  = 3
    7  = 5
    8  = 4
    9  = 0
    10 = 2
    11 = 1
    12 = 0
    13 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
